$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2,6).Value2 = 1424
$ws.Cells.Item(3,6).Value2 = 1406
$ws.Cells.Item(4,6).Value2 = 403
$ws.Cells.Item(5,6).Value2 = 218
$ws.Cells.Item(6,6).Value2 = 663
$ws.Cells.Item(7,6).Value2 = 23
$ws.Cells.Item(8,6).Value2 = 607
$ws.Cells.Item(9,6).Value2 = 471
$ws.Cells.Item(10,6).Value2 = 71
$ws.Cells.Item(11,6).Value2 = 1361
$ws.Cells.Item(12,6).Value2 = 31566
$ws.Cells.Item(12,7).Value2 = 85
$ws.Cells.Item(13,6).Value2 = 6799
$ws.Cells.Item(14,6).Value2 = 101
$ws.Cells.Item(15,4).Value2 = "展贸东路200号 恒达智慧汽车城"
$ws.Cells.Item(15,6).Value2 = 336
$ws.Cells.Item(15,9).Value2 = "//i1.hdslb.com/bfs/openplatform/202407/bC8tPkSW1719820164525.jpeg"
$ws.Cells.Item(16,6).Value2 = 561
$ws.Cells.Item(17,6).Value2 = 322
$ws.Cells.Item(19,6).Value2 = 85
$ws.Cells.Item(20,6).Value2 = 41
$ws.Cells.Item(21,6).Value2 = 425
$ws.Cells.Item(22,6).Value2 = 90
$ws.Cells.Item(23,6).Value2 = 774
$ws.Cells.Item(24,6).Value2 = 310
$ws.Cells.Item(25,6).Value2 = 372
$ws.Cells.Item(26,6).Value2 = 419
$ws.Cells.Item(28,6).Value2 = 178
$ws.Cells.Item(29,6).Value2 = 41
$ws.Cells.Item(30,6).Value2 = 722
$ws.Cells.Item(31,6).Value2 = 278
$ws.Cells.Item(32,6).Value2 = 129
$ws.Cells.Item(33,6).Value2 = 705
$ws.Cells.Item(34,6).Value2 = 102
$ws.Cells.Item(35,3).Value2 = "广州·鸟山明作品《龙珠》40周年only纪念展（取消）"
$ws.Cells.Item(35,7).Value2 = "不可售"
$ws.Cells.Item(36,6).Value2 = 772
$ws.Cells.Item(37,6).Value2 = 276
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2,6).Value2 = 1133
$ws.Cells.Item(5,6).Value2 = 140
$ws.Cells.Item(6,6).Value2 = 288
$ws.Cells.Item(7,6).Value2 = 4308
$ws.Cells.Item(8,6).Value2 = 3
$ws.Cells.Item(9,6).Value2 = 228
$ws.Cells.Item(13,6).Value2 = 38
$ws.Cells.Item(19,6).Value2 = 4279
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2,6).Value2 = 1412
$ws.Cells.Item(3,6).Value2 = 345
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2,6).Value2 = 1412
$ws.Cells.Item(3,6).Value2 = 345
$ws.Cells.Item(4,6).Value2 = 1133
$ws.Cells.Item(5,6).Value2 = 1424
$ws.Cells.Item(6,6).Value2 = 1406
$ws.Cells.Item(7,6).Value2 = 218
$ws.Cells.Item(8,6).Value2 = 663
$ws.Cells.Item(9,6).Value2 = 23
$ws.Cells.Item(10,6).Value2 = 607
$ws.Cells.Item(12,6).Value2 = 71
$ws.Cells.Item(13,6).Value2 = 1361
$ws.Cells.Item(14,6).Value2 = 140
$ws.Cells.Item(15,6).Value2 = 288
$ws.Cells.Item(16,6).Value2 = 228
$ws.Cells.Item(17,6).Value2 = 228
$ws.Cells.Item(20,6).Value2 = 6799
$ws.Cells.Item(21,6).Value2 = 101
$ws.Cells.Item(22,4).Value2 = "展贸东路200号 恒达智慧汽车城"
$ws.Cells.Item(22,6).Value2 = 336
$ws.Cells.Item(22,9).Value2 = "//i1.hdslb.com/bfs/openplatform/202407/bC8tPkSW1719820164525.jpeg"
$ws.Cells.Item(24,6).Value2 = 561
$ws.Cells.Item(25,6).Value2 = 322
$ws.Cells.Item(26,6).Value2 = 38
$ws.Cells.Item(27,6).Value2 = 85
$ws.Cells.Item(29,6).Value2 = 41
$ws.Cells.Item(31,6).Value2 = 425
$ws.Cells.Item(32,6).Value2 = 90
$ws.Cells.Item(33,6).Value2 = 774
$ws.Cells.Item(34,6).Value2 = 310
$ws.Cells.Item(35,6).Value2 = 372
$ws.Cells.Item(36,6).Value2 = 419
$ws.Cells.Item(38,6).Value2 = 178
$ws.Cells.Item(39,6).Value2 = 41
$ws.Cells.Item(40,6).Value2 = 722
$ws.Cells.Item(42,6).Value2 = 278
$ws.Cells.Item(43,6).Value2 = 129
$ws.Cells.Item(44,6).Value2 = 102
$ws.Cells.Item(45,2).NumberFormat = "@"
$ws.Cells.Item(45,2).Value2 = "2024-08-18"
$ws.Cells.Item(45,3).Value2 = "广州·原神×崩坏×绝区零only"
$ws.Cells.Item(45,4).Value2 = "西环路1号 广州岭南会展中心"
$ws.Cells.Item(45,5).Value2 = "2024.08.18 10:00-08.18 17:00"
$ws.Cells.Item(45,6).Value2 = 772
$ws.Cells.Item(45,7).Value2 = 60
$ws.Cells.Item(45,8).Value2 = "https://show.bilibili.com/platform/detail.html?id=87025"
$ws.Cells.Item(45,9).Value2 = "//i0.hdslb.com/bfs/openplatform/202405/lsOq4H701717169339283.png"
$ws.Cells.Item(46,2).NumberFormat = "@"
$ws.Cells.Item(46,2).Value2 = "2024-08-23"
$ws.Cells.Item(46,3).Value2 = "广州·LoveLiveOnly"
$ws.Cells.Item(46,4).Value2 = "芳村大道下市直街1号信义会馆21栋(近白鹅潭风情酒吧街) 信义会馆-21栋"
$ws.Cells.Item(46,5).Value2 = "2024.08.23 10:00-08.23 19:00"
$ws.Cells.Item(46,6).Value2 = 276
$ws.Cells.Item(46,7).Value2 = 68.8
$ws.Cells.Item(46,8).Value2 = "https://show.bilibili.com/platform/detail.html?id=87033"
$ws.Cells.Item(46,9).Value2 = "//i2.hdslb.com/bfs/openplatform/202406/a8shiH411717579829497.jpeg"
$ws.Cells.Item(47,2).NumberFormat = "@"
$ws.Cells.Item(47,2).Value2 = "2024-08-27"
$ws.Cells.Item(47,3).Value2 = "广州·25时主题同人茶会×晓山瑞希生日会"
$ws.Cells.Item(47,4).Value2 = "黄边地铁B出口黄边美食广场1层 胡桃里音乐馆(黄边店)"
$ws.Cells.Item(47,5).Value2 = "2024.08.27 10:00-08.27 16:30"
$ws.Cells.Item(47,6).Value2 = 47
$ws.Cells.Item(47,7).Value2 = 58
$ws.Cells.Item(47,8).Value2 = "https://show.bilibili.com/platform/detail.html?id=87815"
$ws.Cells.Item(47,9).Value2 = "//i1.hdslb.com/bfs/openplatform/202406/rzS5X2Ko1718735908971.png"
$ws.Cells.Item(48,2).NumberFormat = "@"
$ws.Cells.Item(48,2).Value2 = "2024-08-30"
$ws.Cells.Item(48,3).Value2 = "广州·孟京辉经典戏剧作品·黄湘丽主演《一个陌生女人的来信》"
$ws.Cells.Item(48,4).Value2 = "广州市越秀区人民北路696号 广州友谊剧院"
$ws.Cells.Item(48,5).Value2 = "2024.08.30 19:30-08.31 16:30"
$ws.Cells.Item(48,6).Value2 = 12
$ws.Cells.Item(48,7).Value2 = 100
$ws.Cells.Item(48,8).Value2 = "https://show.bilibili.com/platform/detail.html?id=84570"
$ws.Cells.Item(48,9).Value2 = "//i0.hdslb.com/bfs/openplatform/202404/SscDFm1z1713177818070.jpeg"
$ws.Cells.Item(49,2).NumberFormat = "@"
$ws.Cells.Item(49,2).Value2 = "2024-11-05"
$ws.Cells.Item(49,3).Value2 = "广州·2024亚太宠物水族交易会（PSC）国际爬宠展"
$ws.Cells.Item(49,4).Value2 = "新港东路1000号保利世贸博览馆3层 琶洲保利世贸博览馆"
$ws.Cells.Item(49,5).Value2 = "2024.11.05 09:30-11.07 16:30"
$ws.Cells.Item(49,6).Value2 = 21
$ws.Cells.Item(49,7).Value2 = 30
$ws.Cells.Item(49,8).Value2 = "https://show.bilibili.com/platform/detail.html?id=88067"
$ws.Cells.Item(49,9).Value2 = "//i0.hdslb.com/bfs/openplatform/202406/Ej0Rnp201719370264729.jpeg"
$ws.Cells.Item(50,2).NumberFormat = "@"
$ws.Cells.Item(50,2).Value2 = "2024-12-20"
$ws.Cells.Item(50,3).Value2 = "广州·小野丽莎2024“倾爱多彩”唱游世界音乐之旅 纪念专场"
$ws.Cells.Item(50,4).Value2 = "中山纪念堂 中山纪念堂"
$ws.Cells.Item(50,5).Value2 = "2024.12.20 20:00-12.20 22:00"
$ws.Cells.Item(50,6).Value2 = 5
$ws.Cells.Item(50,7).Value2 = 380
$ws.Cells.Item(50,8).Value2 = "https://show.bilibili.com/platform/detail.html?id=87739"
$ws.Cells.Item(50,9).Value2 = "//i0.hdslb.com/bfs/openplatform/202406/HCPstM8c1718868579079.jpeg"
